$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Move Robot21 to location (4, 7) and remove the toolkit."
$ws.Range("B1").Value = "['Robot2']"
$ws.Range("E1").Value = "(4, 7)"

$ws.Range("A2").Value = "Move Robot48 to location (6, 5) and remove the liquid spill."
$ws.Range("E2").Value = "(6, 5)"

$ws.Range("A3").Value = "Move Robot35 to location (3, 3) and remove the large debris."
$ws.Range("E3").Value = "(3, 3)"

$ws.Range("A4").Value = "Move Robot15 to location (3, 4) and remove the dust."
$ws.Range("E4").Value = "(3, 4)"

$ws.Range("A5").Value = "Move Robot29 to location (9, 10) and remove the grass."
$ws.Range("E5").Value = "(9, 10)"

$ws.Range("A6").Value = "Move Robot31 to location (8, 12) and remove the small debris."
$ws.Range("B6").Value = "['Robot8', 'Robot50']"
$ws.Range("E6").Value = "(8, 12)"

$ws.Range("A7").Value = "Move Robot13 to location (7, 5) and remove the vehicle."
$ws.Range("E7").Value = "(7, 5)"

$ws.Range("A8").Value = "Move Robot50 to location (5, 12) and remove the construction materials."
$ws.Range("B8").Value = "['Robot22', 'Robot9', 'Robot13']"
$ws.Range("E8").Value = "(5, 12)"

$ws.Range("A9").Value = "Move Robot9 to location (11, 2) and remove the tree branches."
$ws.Range("B9").Value = "['Robot24']"
$ws.Range("E9").Value = "(11, 2)"

$ws.Range("A10").Value = "Move Robot40 to location (10, 3) and remove the screws."
$ws.Range("E10").Value = "(10, 3)"
